# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest scraped counts.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 7952
    $ws.Range("F10").Value = 178
    $ws.Range("F11").Value = 235
    $ws.Range("F14").Value = 1975
    $ws.Range("F20").Value = 20
}
